# Add institution-import support (#1)
#  - rename the original sheet to "העברות" (Transfers)
#  - add a new "מוסד" (Institution) sheet with a header row + one sample
#    row: institution id / sending-institution id / name / identifying
#    number
#  - keep the new sheet's "identifying number" column formatted as text
#    (numbers like "001" must not lose their leading zero)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: "Sheet1" -> "העברות" --------------------------------------
$ws1.Name = "העברות"

# Column width tweaks observed on the transfers sheet (closest values the
# ColumnWidth API can hit, which snaps to 1/6-character increments)
$ws1.Columns.Item(1).ColumnWidth = 9.333333333333332
$ws1.Columns.Item(5).ColumnWidth = 10.333333333333332

$ws1.Range("E10").Select() | Out-Null

# --- New sheet: "מוסד" ---------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "מוסד"

$ws2.Range("A1").Value = "מזהה מוסד"
$ws2.Range("B1").Value = "מזהה מוסד שולח"
$ws2.Range("C1").Value = "שם מוסד"

$ws2.Range("D1:D2").NumberFormat = "@"
$ws2.Range("D1").Value = "מספר מזהה"

$ws2.Range("A2").Value = 12345678
$ws2.Range("B2").Value = 12345
$ws2.Range("D2").Value = "001"
$ws2.Range("C2").Value = "החברה שלי בע""מ"

$ws2.Columns.Item(1).ColumnWidth = 21.0
$ws2.Columns.Item(2).ColumnWidth = 14.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(4).ColumnWidth = 10.833333333333332

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("B5").Select() | Out-Null

$ws1.Select() | Out-Null
